# PO Clean up - updated the existing test case data on the
# InventoryRequisition sheet (QuantityBefore / QuantityAfter columns
# for the two sample rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryRequisition")

# Columns H/I store these as text (e.g. "1559.0"), not numbers, so force
# text formatting before/after the write to avoid Excel auto-converting
# the value to a number and to keep the cell's original (default) style.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "1639.0"
$ws.Range("H2").Style = "Normal"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1641.0"
$ws.Range("I2").Style = "Normal"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "955.0"
$ws.Range("H3").Style = "Normal"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "957.0"
$ws.Range("I3").Style = "Normal"
